# Refresh crypto symbol list prices/volumes (Price + Volume(1h) columns)
# to match the "Updated symbol list ... with GitHub Actions" commit.
# Values are stored as literal text (e.g. "258.48", "0.73%") in the source
# workbook, so each assignment uses a leading apostrophe to force Excel to
# keep the text verbatim instead of re-parsing it as a number/percentage,
# then resets the cell style back to Normal so no stray "quote prefix"
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '258.48'
    'E2' = '0.73%'
    'E3' = '-0.82%'
    'D4' = '4.645'
    'E4' = '1.92%'
    'D5' = '0.05979'
    'E5' = '1.28%'
    'D6' = '6.643'
    'E6' = '0.29%'
    'D7' = '0.8565'
    'E7' = '-0.39%'
    'D8' = '0.9215'
    'E8' = '-0.96%'
    'D9' = '0.1387'
    'E9' = '-1.49%'
    'D10' = '0.04547'
    'E10' = '25.35%'
    'D11' = '0.07006'
    'E11' = '-0.90%'
    'D12' = '0.03051'
    'E12' = '-5.51%'
    'D13' = '0.09117'
    'E13' = '-0.96%'
    'D14' = '0.001533'
    'E14' = '-1.66%'
    'D15' = '0.0006062'
    'E15' = '0.04%'
    'D16' = '0.006185'
    'E16' = '1.62%'
    'D17' = '3.445'
    'E17' = '-1.99%'
    'D18' = '3.151'
    'E18' = '-1.38%'
    'E19' = '-2.19%'
    'D20' = '0.3106'
    'E20' = '1.63%'
    'E21' = '0.87%'
    'D22' = '4.026'
    'E22' = '4.58%'
    'D23' = '0.04227'
    'E23' = '0.16%'
    'E24' = '-0.24%'
    'D25' = '0.004028'
    'E25' = '-5.88%'
    'D26' = '0.0001200'
    'E26' = '-0.07%'
    'E27' = '-11.62%'
    'D40' = '0.03831'
    'E40' = '0.20%'
    'E41' = '1.04%'
    'D42' = '0.003761'
    'E42' = '-4.80%'
    'E43' = '0.34%'
    'D44' = '0.01504'
    'E44' = '32.80%'
    'D45' = '0.00005110'
    'E45' = '-5.98%'
    'E46' = '-0.02%'
    'E47' = '-17.01%'
    'D48' = '0.1576'
    'E48' = '58.42%'
    'E49' = '-0.02%'
    'E50' = '-0.02%'
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $updates[$cellRef]
    $cell.Style = "Normal"
}
